# Generate Report for Handoff
# Updates the status from "In Translation" to "Ready for handoff" and
# refreshes the related "Latest ... Datetime" timestamps on the Overview,
# zh-cn and de-de sheets, then lets the column widths auto-fit to the
# new (longer) status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-11-29 05:01:56"

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-11-29 05:01:43"

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Widen the columns that now hold the longer "Ready for handoff" text, to
# match what Excel's own AutoFit would compute for the new status string.
$wsOverview.Range("E:E").ColumnWidth = 16.33
$wsOverview.Range("F:F").ColumnWidth = 16.33
$wsZhCn.Range("C:C").ColumnWidth = 16.33
$wsDeDe.Range("C:C").ColumnWidth = 16.33
